# Update Name of Algo
# Update specific values in column B of Sheet1 to reflect the revised
# KNN imputation results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value = 6.345999999999999
$ws.Range("B18").Value = 5.126
$ws.Range("B20").Value = 6.967000000000001
$ws.Range("B27").Value = 6.161
$ws.Range("B69").Value = 5.306999999999999
$ws.Range("B76").Value = 6.308
$ws.Range("B82").Value = 5.366000000000001
